$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("research_bar")

# Update the "Undergraduate Research Assistant" entries so the 'with' / 'where'
# columns read "Department of Chemistry" / " University of Puerto Rico, San Juan, PR"
# instead of the previous split text (modeled after L. Abad's CV).
$ws.Range("G6").Value = "Department of Chemistry"
$ws.Range("H6").Value = " University of Puerto Rico, San Juan, PR"
$ws.Range("G10").Value = "Department of Chemistry"
$ws.Range("H10").Value = " University of Puerto Rico, San Juan, PR"

# Widen column H to fit the new text and select it as the active cell.
# (46.7109375 characters is the author's target; the COM width setter here
# snaps to 1/6-character increments, so 45.8333... is the closest input that
# lands on the nearest achievable grid value, 46.6666...)
$ws.Columns.Item(8).ColumnWidth = 45.8333333333333
$ws.Range("H18").Select()
